# Update STEO snapshot from "January 2017" to "February 2017"
# and refresh the underlying history/forecast data (rows 39-41, 158-183).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig24")

# --- Title / source text strings -------------------------------------------------
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A184").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Annual average price / growth table (F/G columns) ---------------------------
$ws.Range("F39").Value = 12.547745127000001
$ws.Range("F40").Value = 12.927405175000001
$ws.Range("F41").Value = 13.240067381999999

# --- Monthly history (B) / forecast (C) series ------------------------------------
$ws.Range("B158").Value = 12.75
$ws.Range("B159").Value = 12.215070000000001
$ws.Range("C159").Value = "#N/A"

$ws.Range("B160").Value = 12.25662
$ws.Range("C160").Value = 12.25662

$ws.Range("C161").Value = 12.55918
$ws.Range("C162").Value = 12.73115
$ws.Range("C163").Value = 12.63547
$ws.Range("C164").Value = 13.01451
$ws.Range("C165").Value = 13.043049999999999
$ws.Range("C166").Value = 13.106249999999999
$ws.Range("C167").Value = 13.35501
$ws.Range("C168").Value = 13.433949999999999
$ws.Range("C169").Value = 13.00178
$ws.Range("C170").Value = 13.2156
$ws.Range("C171").Value = 12.64715
$ws.Range("C172").Value = 12.76178
$ws.Range("C173").Value = 13.05621
$ws.Range("C174").Value = 13.169790000000001
$ws.Range("C175").Value = 13.00924
$ws.Range("C176").Value = 13.3367
$ws.Range("C177").Value = 13.309559999999999
$ws.Range("C178").Value = 13.3391
$ws.Range("C179").Value = 13.57114
$ws.Range("C180").Value = 13.64565
$ws.Range("C181").Value = 13.23039
$ws.Range("C182").Value = 13.47226
$ws.Range("C183").Value = 12.93544
